$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for a new "Metryki architektoniczne" section title + spacer
#    row above the existing (Testowalnosc / Konwencje / n.d) table. The
#    original header row (row 3) and its data rows (4-6) shift down by one
#    (-> 4-7).
# ---------------------------------------------------------------------------
$ws.Rows("3:3").Insert() | Out-Null

# Section title in row 2.
$titleCell = $ws.Range("A2")
$titleCell.Font.Bold = $true
$titleCell.Font.Size = 14
$titleCell.Value = "Metryki architektoniczne"

# Small spacer row (row 3) just under the title - bigger default font,
# bold on the first cell, a bit taller than normal.
$ws.Rows("3:3").RowHeight = 15.75
$spacerCell = $ws.Range("A3")
$spacerCell.Font.Size = 12
$spacerCell.Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. New "Metryki biznesowe" section: title + a copy of the same table
#    header + three new data rows.
# ---------------------------------------------------------------------------
$biz = $ws.Range("A11")
$biz.Font.Bold = $true
$biz.Font.Size = 14
$biz.Value = "Metryki biznesowe"

# Re-use the existing header formatting/content (row 4: Typ drivera / Driver
# / Jak rozumiemy? / Wartosc aktualna / Limit / Cel / Idea») for the new
# table header on row 13.
$ws.Range("B4:H4").Copy($ws.Range("B13:H13")) | Out-Null

# Row 14 - "Czas dodania wpisu dziennika"
$ws.Range("A14").Value = "(metryka dlużna)"
$ws.Range("A14").Value = "(metryka dłużna)"
$ws.Range("B14").Value = "Atrybut jakościowy"
$ws.Range("C14").WrapText = $true
$ws.Range("C14").Value = "Czas dodania wpisu dziennika"
$ws.Range("E14").Style = "Neutral"
$ws.Range("E14").Value = "?"
$ws.Range("D14").Value = "Czas w sekunadach potrzebny na dodanie nowego wpisu dziennika"
$ws.Range("F14").Value = "?"
$ws.Range("G14").Value = "?"
$ws.Range("H14").Value = "?"

# Row 15 - "Czas dodania klienta wraz z pojazdem"
$ws.Range("A15").Value = "(metryka dłużna)"
$ws.Range("B15").Value = "Atrybut jakościowy"
$ws.Range("C15").WrapText = $true
$ws.Range("C15").Value = "Czas dodania klienta wraz z pojazdem"

# Row 16 - "Czas dodania standardowej naprawy"
$ws.Range("A16").Value = "(metryka dłużna)"
$ws.Range("B16").Value = "Atrybut jakościowy"
$ws.Range("C16").WrapText = $true
$ws.Range("C16").Value = "Czas dodania standardowej naprawy"

$ws.Range("D15").Value = "Czas w sekunadach potrzebny na dodanie nowego klienta wraz z pojazdem"
$ws.Range("D16").Value = "Czas w sekunadach potrzebny na dodanie standardowej naprawy"

$ws.Range("E15").Style = "Neutral"
$ws.Range("E15").Value = "?"
$ws.Range("F15").Value = "?"
$ws.Range("G15").Value = "?"
$ws.Range("H15").Value = "?"

$ws.Range("E16").Style = "Neutral"
$ws.Range("E16").Value = "?"
$ws.Range("F16").Value = "?"
$ws.Range("G16").Value = "?"
$ws.Range("H16").Value = "?"

$ws.Range("I15").Value = "Metryka weryfikująca sukces"

$ws.Range("I16").Select() | Out-Null
